{"js": "// Replace each three-digit-by-one-digit multiplication expression\n// with its new value, preserving the run formatting (font/size).\nconst replacements = [\n  ['614\u00d77=', '282\u00d73='],\n  ['721\u00d79=', '296\u00d77='],\n  ['731\u00d72=', '569\u00d77='],\n  ['384\u00d77=', '867\u00d79='],\n  ['332\u00d73=', '378\u00d72='],\n  ['790\u00d78=', '478\u00d72='],\n  ['280\u00d74=', '923\u00d72='],\n  ['384\u00d72=', '825\u00d73='],\n  ['370\u00d73=', '710\u00d72='],\n  ['761\u00d78=', '145\u00d78='],\n  ['743\u00d78=', '525\u00d78='],\n  ['208\u00d78=', '954\u00d79='],\n  ['140\u00d74=', '710\u00d79='],\n  ['179\u00d73=', '905\u00d76='],\n  ['812\u00d79=', '434\u00d76='],\n  ['166\u00d72=', '224\u00d72='],\n  ['614\u00d73=', '936\u00d74='],\n  ['330\u00d75=', '784\u00d76='],\n  ['838\u00d78=', '668\u00d77='],\n  ['374\u00d77=', '399\u00d72='],\n  ['266\u00d79=', '900\u00d74='],\n  ['804\u00d76=', '640\u00d72='],\n  ['545\u00d77=', '356\u00d79='],\n  ['397\u00d75=', '250\u00d75='],\n  ['973\u00d73=', '563\u00d72='],\n];\n\nfor (const [findText, replaceText] of replacements) {\n  const results = context.document.body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication expression\n# with its new value, preserving run formatting (font/size).\n$d = $word.ActiveDocument\n$mult = [char]0x00D7\n$replacements = @(\n    @(\"614${mult}7=\", \"282${mult}3=\"),\n    @(\"721${mult}9=\", \"296${mult}7=\"),\n    @(\"731${mult}2=\", \"569${mult}7=\"),\n    @(\"384${mult}7=\", \"867${mult}9=\"),\n    @(\"332${mult}3=\", \"378${mult}2=\"),\n    @(\"790${mult}8=\", \"478${mult}2=\"),\n    @(\"280${mult}4=\", \"923${mult}2=\"),\n    @(\"384${mult}2=\", \"825${mult}3=\"),\n    @(\"370${mult}3=\", \"710${mult}2=\"),\n    @(\"761${mult}8=\", \"145${mult}8=\"),\n    @(\"743${mult}8=\", \"525${mult}8=\"),\n    @(\"208${mult}8=\", \"954${mult}9=\"),\n    @(\"140${mult}4=\", \"710${mult}9=\"),\n    @(\"179${mult}3=\", \"905${mult}6=\"),\n    @(\"812${mult}9=\", \"434${mult}6=\"),\n    @(\"166${mult}2=\", \"224${mult}2=\"),\n    @(\"614${mult}3=\", \"936${mult}4=\"),\n    @(\"330${mult}5=\", \"784${mult}6=\"),\n    @(\"838${mult}8=\", \"668${mult}7=\"),\n    @(\"374${mult}7=\", \"399${mult}2=\"),\n    @(\"266${mult}9=\", \"900${mult}4=\"),\n    @(\"804${mult}6=\", \"640${mult}2=\"),\n    @(\"545${mult}7=\", \"356${mult}9=\"),\n    @(\"397${mult}5=\", \"250${mult}5=\"),\n    @(\"973${mult}3=\", \"563${mult}2=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n"}
